$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 699 (shifts existing rows 699-780 down to 700-781)
$ws.Rows.Item(699).Insert()

# Populate the newly inserted row 699 with the new record's data
$ws.Cells.Item(699, 1).Value = 6
$ws.Cells.Item(699, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(699, 3).Value = "Metropolitana"
$ws.Cells.Item(699, 4).Value = 45212
$ws.Cells.Item(699, 5).Value = 13
$ws.Cells.Item(699, 6).Value = 100112039
$ws.Cells.Item(699, 7).Value = "Ciboulette"
$ws.Cells.Item(699, 8).Value = "Sin especificar"
$ws.Cells.Item(699, 9).Value = "Primera"
$ws.Cells.Item(699, 10).Value = 650
$ws.Cells.Item(699, 11).Value = 900
$ws.Cells.Item(699, 12).Value = 1000
$ws.Cells.Item(699, 13).Value = 954
$ws.Cells.Item(699, 14).Value = "$/docena de atados"
$ws.Cells.Item(699, 15).Value = "Región Metropolitana"
$ws.Cells.Item(699, 16).Value = 318
$ws.Cells.Item(699, 17).Value = 3
$ws.Cells.Item(699, 18).Value = "Hortaliza"
